$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-05-07 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-05-08 Thursday", 2) | Out-Null

# Update the division problems table, cell by cell (row, col are 1-based)
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="60÷4="},
    @{Row=1;  Col=2; Text="29÷3="},
    @{Row=1;  Col=3; Text="76÷8="},
    @{Row=1;  Col=4; Text="84÷3="},
    @{Row=1;  Col=5; Text="64÷2="},

    @{Row=5;  Col=1; Text="24÷8="},
    @{Row=5;  Col=2; Text="10÷8="},
    @{Row=5;  Col=3; Text="76÷3="},
    @{Row=5;  Col=4; Text="93÷5="},
    @{Row=5;  Col=5; Text="37÷7="},

    @{Row=9;  Col=1; Text="30÷9="},
    @{Row=9;  Col=2; Text="46÷5="},
    @{Row=9;  Col=3; Text="46÷3="},
    @{Row=9;  Col=4; Text="50÷9="},
    @{Row=9;  Col=5; Text="76÷9="},

    @{Row=13; Col=1; Text="61÷5="},
    @{Row=13; Col=2; Text="14÷2="},
    @{Row=13; Col=3; Text="53÷8="},
    @{Row=13; Col=4; Text="31÷3="},
    @{Row=13; Col=5; Text="39÷9="},

    @{Row=17; Col=1; Text="13÷2="},
    @{Row=17; Col=2; Text="28÷9="},
    @{Row=17; Col=3; Text="23÷4="},
    @{Row=17; Col=4; Text="23÷5="},
    @{Row=17; Col=5; Text="59÷8="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
